# Generate Report for Handback
# Update the timestamps recorded for the handoff/handback report.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file.
$wsOverview.Range("G2").Value = "2016-08-20 07:07:58"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime".
$wsZhCn.Range("H2").Value = "2016-08-20 07:07:55"
$wsZhCn.Range("K2").Value = "2016-08-20 07:08:14"

# de-de sheet: "Correspond Handback DateTime".
$wsDeDe.Range("K2").Value = "2016-08-20 07:08:20"
